# Natmi LR-pairs output update (Slurp1-Chrna7), following Dr Hou's advice:
# the single FAPs->FAPs/sCs summary rows are expanded into the full
# per-sending/target-cluster combination (FAPs and sCs crossed with
# FAPs/sCs/ECs), with refreshed NATMI statistics for each combination.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A..D = Sending cluster / Ligand symbol / Receptor symbol / Target cluster
#          E..T = the 16 NATMI numeric statistics
$rows = @(
    @("FAPs","Slurp1","Chrna7","ECs", 3,1,0.9671793333333333,2.901538,0.9567964966834895,0.9567964966834897,1,0.3333333333333333,0.05437633333333333,0.163129,0.01289405845289636,0.01289405845289636,0.05259166582244444,0.473324992402,0.01233698995576337,0.01233698995576338),
    @("FAPs","Slurp1","Chrna7","FAPs", 3,1,0.9671793333333333,2.901538,0.9567964966834895,0.9567964966834897,3,1,0.2159606666666667,0.647882,0.05120995272808269,0.05120995272808269,0.2088726936128889,1.879854242516,0.04899750336555663,0.04899750336555663),
    @("FAPs","Slurp1","Chrna7","sCs", 3,1,0.9671793333333333,2.901538,0.9567964966834895,0.9567964966834897,3,1,3.946825,11.840475,0.9358959888190209,0.9358959888190209,3.817287572283333,34.35558815055,0.8954620033621695,0.8954620033621696),
    @("sCs","Slurp1","Chrna7","ECs", 1,0.3333333333333333,0.04367233333333333,0.131017,0.04320350331651034,0.04320350331651034,1,0.3333333333333333,0.05437633333333333,0.163129,0.01289405845289636,0.01289405845289636,0.002374741354777778,0.021372672193,0.0005570684971329861,0.0005570684971329861),
    @("sCs","Slurp1","Chrna7","FAPs", 1,0.3333333333333333,0.04367233333333333,0.131017,0.04320350331651034,0.04320350331651034,3,1,0.2159606666666667,0.647882,0.05120995272808269,0.05120995272808269,0.009431506221555556,0.08488355599399999,0.002212449362526058,0.002212449362526058),
    @("sCs","Slurp1","Chrna7","sCs", 1,0.3333333333333333,0.04367233333333333,0.131017,0.04320350331651034,0.04320350331651034,3,1,3.946825,11.840475,0.9358959888190209,0.9358959888190209,0.1723670570083333,1.551303513075,0.04043398545685129,0.04043398545685129)
)

$rowIndex = 2
foreach ($rowData in $rows) {
    $colIndex = 1
    foreach ($value in $rowData) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $value
        $colIndex++
    }
    $rowIndex++
}
